$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '28.986.76'
Set-TextValue 'E2' '  -3.88%  '

Set-TextValue 'D3' '1.958.45'
Set-TextValue 'E3' '  -5.70%  '

Set-TextValue 'E4' '  -0.20%  '

Set-TextValue 'D5' '326.44'
Set-TextValue 'E5' '  -3.18%  '

Set-TextValue 'E6' '  +0.03%  '

Set-TextValue 'D7' '0.4953'
Set-TextValue 'E7' '  -5.20%  '

Set-TextValue 'D8' '0.4192'
Set-TextValue 'E8' '  -3.51%  '

Set-TextValue 'D9' '52.88'
Set-TextValue 'E9' '  -3.83%  '

Set-TextValue 'D10' '0.09236'
Set-TextValue 'E10' '  -0.71%  '

Set-TextValue 'D11' '1.094'
Set-TextValue 'E11' '  -6.10%  '

Set-TextValue 'D12' '22.70'
Set-TextValue 'E12' '  -6.69%  '

Set-TextValue 'D13' '1.959.40'
Set-TextValue 'E13' '  -4.46%  '

Set-TextValue 'D14' '6.430'
Set-TextValue 'E14' '  -5.71%  '

Set-TextValue 'D15' '7.814'
Set-TextValue 'E15' '  -6.98%  '

Set-TextValue 'E16' '  +0.20%  '

Set-TextValue 'B17' 'ShibaInu'
Set-TextValue 'C17' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D17' '0.00001096'
Set-TextValue 'E17' '  -4.71%  '

Set-TextValue 'B18' 'Litecoin'
Set-TextValue 'C18' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D18' '91.07'
Set-TextValue 'E18' '  -9.36%  '

Set-TextValue 'D19' '0.06682'
Set-TextValue 'E19' '  -0.27%  '

Set-TextValue 'D20' '19.16'
Set-TextValue 'E20' '  -7.51%  '

Set-TextValue 'E21' '  +0.10%  '

Set-TextValue 'D22' '5.924'
Set-TextValue 'E22' '  -5.18%  '

Set-TextValue 'D23' '29.029.19'
Set-TextValue 'E23' '  -3.88%  '

Set-TextValue 'D24' '11.95'
Set-TextValue 'E24' '  -2.71%  '

Set-TextValue 'E25' '  -2.16%  '

Set-TextValue 'D26' '2.200.84'
Set-TextValue 'E26' '  -4.01%  '

Set-TextValue 'D27' '20.51'
Set-TextValue 'E27' '  -5.24%  '

Set-TextValue 'D28' '155.13'
Set-TextValue 'E28' '  -4.26%  '

Set-TextValue 'D29' '6.263'
Set-TextValue 'E29' '  -7.35%  '

Set-TextValue 'D30' '2.234'
Set-TextValue 'E30' '  -8.99%  '

Set-TextValue 'D31' '125.93'
Set-TextValue 'E31' '  -5.15%  '

Set-TextValue 'D32' '1.036'
Set-TextValue 'E32' '  -7.16%  '

Set-TextValue 'D33' '0.09804'
Set-TextValue 'E33' '  -5.74%  '

Set-TextValue 'D34' '1.508'
Set-TextValue 'E34' '  -7.06%  '

Set-TextValue 'D35' '5.803'
Set-TextValue 'E35' '  -6.58%  '

Set-TextValue 'D36' '3.674'
Set-TextValue 'E36' '  -5.77%  '

Set-TextValue 'D37' '0.02415'
Set-TextValue 'E37' '  -6.62%  '

Set-TextValue 'D38' '1.314'
Set-TextValue 'E38' '  -0.98%  '

Set-TextValue 'B39' 'Hedera'
Set-TextValue 'C39' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D39' '0.06326'
Set-TextValue 'E39' '  -5.05%  '

Set-TextValue 'B40' 'FraxShare'
Set-TextValue 'C40' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D40' '8.945'
Set-TextValue 'E40' '  -8.39%  '

Set-TextValue 'D41' '0.6409'
Set-TextValue 'E41' '  -6.72%  '

Set-TextValue 'D42' '11.33'
Set-TextValue 'E42' '  -8.33%  '

Set-TextValue 'D43' '0.1970'
Set-TextValue 'E43' '  -9.38%  '

Set-TextValue 'E44' '  +0.09%  '

Set-TextValue 'D45' '1.363'
Set-TextValue 'E45' '  +4.03%  '

Set-TextValue 'D46' '0.6176'
Set-TextValue 'E46' '  -7.27%  '

Set-TextValue 'D47' '13.29'
Set-TextValue 'E47' '  -6.53%  '

Set-TextValue 'D48' '2.183'
Set-TextValue 'E48' '  -6.25%  '

Set-TextValue 'D49' '3.457'
Set-TextValue 'E49' '  -4.44%  '

Set-TextValue 'D50' '0.00000000334'
Set-TextValue 'E50' '  -3.63%  '

Set-TextValue 'D51' '0.06994'
Set-TextValue 'E51' '  -2.87%  '

Write-Host "Done updating cryptos sheet"